# Updated cryptos list on Mon Feb 20 05:49:01 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) for each coin row.
# Column D values are textual (dotted thousands separators, e.g. "24.465.04")
# so for entries that look like plain decimals we force NumberFormat "@"
# (Text) first to stop Excel from auto-coercing the assigned string into a
# numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.465.04"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "1.686.27"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.87"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3913"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4030"
$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("E9").Value = "  -2.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.38"
$ws.Range("E11").Value = "  -1.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08775"
$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.50"
$ws.Range("E13").Value = "  +11.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.454"
$ws.Range("E14").Value = "  +1.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.155"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("D17").Value = "1.681.33"
$ws.Range("E17").Value = "  -2.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "98.10"
$ws.Range("E18").Value = "  -2.47%  "

$ws.Range("E19").Value = "  +2.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.28"
$ws.Range("E20").Value = "  +2.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.285"
$ws.Range("E21").Value = "  +2.73%  "

$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.31"
$ws.Range("E23").Value = "  -1.58%  "

$ws.Range("D24").Value = "24.452.96"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.023"
$ws.Range("E25").Value = "  -6.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.338"
$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.63"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.444"
$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.343"
$ws.Range("E30").Value = "  +3.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "138.23"
$ws.Range("E31").Value = "  +1.00%  "

$ws.Range("D32").Value = "1.864.38"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08735"
$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.254"
$ws.Range("E34").Value = "  -4.16%  "

$ws.Range("E35").Value = "  -3.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.085"
$ws.Range("E36").Value = "  +5.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.03004"
$ws.Range("E37").Value = "  +8.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2791"
$ws.Range("E38").Value = "  +1.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.85"
$ws.Range("E39").Value = "  -3.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09141"
$ws.Range("E40").Value = "  -0.87%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8029"
$ws.Range("E41").Value = "  +3.87%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.15"
$ws.Range("E42").Value = "  -3.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.469"
$ws.Range("E43").Value = "  +0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.50"
$ws.Range("E44").Value = "  +9.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.650"
$ws.Range("E45").Value = "  +2.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7249"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.262"
$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("E48").Value = "  +6.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9989"
$ws.Range("E49").Value = "  -0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.35"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("E51").Value = "  +0.98%  "
